# Apply the edits described by the diff:
#  1. Add a new shared string "verifyloanproduct" and use it as the value of
#     A1 on the ProductLoanOutput sheet (replacing "productname").
#  2. On ProductLoanInput, restyle B7 ("currency" value "US Dollar") to use
#     the plain/general style (same style as e.g. B15), instead of the
#     numeric style it incorrectly had.
#  3. Update the selection/view state on ProductLoanInput to focus A7:B7
#     (instead of A27), with no special top-left scroll position.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# --- ProductLoanOutput: A1 becomes "verifyloanproduct" --------------------
$ws2.Range("A1").Value = "verifyloanproduct"

# --- ProductLoanInput: fix formatting of B7 (currency value) --------------
# Copy the format from B15, which already uses the desired plain style,
# onto B7, without touching B7's value.
$ws1.Range("B15").Copy()
$ws1.Range("B7").PasteSpecial(-4122)   # xlPasteFormats

# --- ProductLoanInput: update the active selection / view -----------------
$ws2.Activate()
$ws2.Range("A1").Select()

$ws1.Activate()
$ws1.Range("A7:B7").Select()

Write-Host "edits applied"
